# worked on evaluatie gj2018
$wb = $excel.ActiveWorkbook

$wsReal = $wb.Worksheets.Item("realisatie")
$wsProg = $wb.Worksheets.Item("prognose")

# Fill in the new 2018 realisation row (row 20) that used to be an empty
# placeholder cell (B20) on the "realisatie" sheet.
$wsReal.Range("B20").ClearFormats()
$wsReal.Cells.Item(20, 1).Value = 2018     # A20
$wsReal.Cells.Item(20, 2).Value = 497455   # B20
$wsReal.Cells.Item(20, 3).Value = 9200     # C20
$wsReal.Cells.Item(20, 4).Value = 150000   # D20
$wsReal.Cells.Item(20, 5).Value = 2500     # E20
$wsReal.Cells.Item(20, 6).Value = 19500    # F20

# Reflect the new working selection/active sheet: the user ended up back on
# "realisatie" with G20 selected, while "prognose" was left selected at E9.
$wsProg.Range("E9").Select()
$wsReal.Select()
$wsReal.Range("G20").Select()
